$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a genuine run boundary at the edges of a [start, end) range by
# toggling a character property on and back off. Because of how this runtime
# reconciles adjacent runs that share identical formatting, simply writing
# text leaves every edited/neighbouring run with matching rPr merged back
# into one <w:r>. Touching (and then restoring) a direct-formatting property
# on the exact sub-range keeps that sub-range as its own run without leaving
# any visible formatting residue.
# ---------------------------------------------------------------------------
function Mark-RunBoundary($rangeStart, $rangeEnd) {
    $seg = $d.Range($rangeStart, $rangeEnd)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# 1) "Versão: 2.0" -> "Versão: 3.0", with "3" and ".0" kept as distinct runs
# ---------------------------------------------------------------------------
$verRange = $d.Content
$verRange.Find.Execute("2.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verStart = $verRange.Start
$verRange.Text = "3.0"

$verMid = $verStart + 1
$verEnd = $verStart + 3
Mark-RunBoundary $verStart $verMid  # "3"
Mark-RunBoundary $verMid $verEnd    # ".0"

# ---------------------------------------------------------------------------
# 2) "Revisão: 03/08/2023" -> "Revisão: 12/04/2025", split across five runs:
#    "12" / "/0" / "4" / "/202" / "5"
# ---------------------------------------------------------------------------
$dateRange = $d.Content
$dateRange.Find.Execute("03/08/2023", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dateStart = $dateRange.Start
$dateRange.Text = "12/04/2025"

$segmentLengths = @(2, 2, 1, 4, 1)
$pos = $dateStart
foreach ($len in $segmentLengths) {
    $segEnd = $pos + $len
    Mark-RunBoundary $pos $segEnd
    $pos = $segEnd
}

# ---------------------------------------------------------------------------
# 3) Merge the four runs spelling "RF" + "-02" + "2" + " Configurar conexão
#    remota" back into a single run "RF-022 Configurar conexão remota"
# ---------------------------------------------------------------------------
$rfRange = $d.Content
$rfRange.Find.Execute("RF-022 Configurar conexão remota", $true, $false, $false, $false, $false, $true, 1, $false, "RF-022 Configurar conexão remota", 2)
